# Insert a new weekly price record into the "Zanahoria" (carrot) price
# sheet. This pushes the existing rows 399..469 down by one (to 400..470)
# and fills the freshly-opened row 399 with a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 399 downward to make room for the new record.
$ws.Rows.Item(399).Insert()

# Populate the new row with the new weekly observation.
$ws.Cells.Item(399, 1).Value = 8
$ws.Cells.Item(399, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(399, 3).Value = "Coquimbo"
$ws.Cells.Item(399, 4).Value = 44951
$ws.Cells.Item(399, 5).Value = 4
$ws.Cells.Item(399, 6).Value = 100114013
$ws.Cells.Item(399, 7).Value = "Zanahoria"
$ws.Cells.Item(399, 8).Value = "Sin especificar"
$ws.Cells.Item(399, 9).Value = "Primera"
$ws.Cells.Item(399, 10).Value = 600
$ws.Cells.Item(399, 11).Value = 5800
$ws.Cells.Item(399, 12).Value = 6000
$ws.Cells.Item(399, 13).Value = 5900
$ws.Cells.Item(399, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(399, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(399, 16).Value = 295
$ws.Cells.Item(399, 17).Value = 20
$ws.Cells.Item(399, 18).Value = "Hortaliza"
